# Updated symbol list (price refresh + WazirX.. shuffle) per GitHub Actions run.
# Numeric-looking values in column D are entered with a leading apostrophe so
# Excel keeps them as text (matching the workbook's existing text-stored
# price column) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'240.78"
$ws.Range("D3").Value = "'21.39"
$ws.Range("D4").Value = "'5.132"
$ws.Range("D5").Value = "'0.05535"
$ws.Range("D6").Value = "'3.374"
$ws.Range("D7").Value = "'6.372"
$ws.Range("D8").Value = "'0.8049"
$ws.Range("D9").Value = "'0.9422"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1385"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07243"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D12").Value = "'0.03065"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03075"
$ws.Range("E13").Value = "12BitrueCoinBTR"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09317"
$ws.Range("E14").Value = "13BitMartTokenBMX"
$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D15").Value = "'3.614"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001620"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "'0.04709"
$ws.Range("E17").Value = "16CoinExTokenCET"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0005766"
$ws.Range("E18").Value = "17OneONE"
$ws.Range("D19").Value = "'0.006396"
$ws.Range("D21").Value = "'0.001046"
$ws.Range("D22").Value = "'0.0001506"
$ws.Range("D23").Value = "'0.0003115"
$ws.Range("D24").Value = "'3.751"
$ws.Range("D25").Value = "'2.102"
$ws.Range("D26").Value = "'0.3253"
$ws.Range("D27").Value = "'0.1292"
$ws.Range("D40").Value = "'0.03868"
$ws.Range("D41").Value = "'0.006908"
$ws.Range("D42").Value = "'0.1028"
$ws.Range("D43").Value = "'0.003097"
$ws.Range("D44").Value = "'0.008274"
$ws.Range("D45").Value = "'0.00005967"
$ws.Range("D46").Value = "'0.00000000754"
$ws.Range("D47").Value = "'0.0005525"
$ws.Range("D48").Value = "'0.6858"
$ws.Range("D49").Value = "'0.1039"
$ws.Range("E49").Value = "48BOLOBOLOBestin24h"
$ws.Range("D50").Value = "'0.00002110"
$ws.Range("D51").Value = "'0.01015"
